$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.339.92'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.809.40'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'702.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'171.48"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').Value = '3.809.89'
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = "'7.50"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('D12').Value = "'0.477"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.27%  '
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('D14').Value = "'36.04"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '4.449.89'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '3.814.47'
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('D17').Value = '71.295.76'
$ws.Range('D18').Value = "'7.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'17.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = "'0.115"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = "'514.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.22%  '
$ws.Range('D22').Value = "'10.50"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').Value = "'0.714"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').Value = "'84.08"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('D25').Value = "'0.0000143"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.87%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').Value = '3.956.91'
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').Value = "'10.35"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = "'2.03"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.84%  '
$ws.Range('D31').Value = "'3.01"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.45%  '
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('D34').Value = "'29.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').Value = "'0.174"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('D36').Value = "'9.15"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '3.770.56'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').Value = "'6.50"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.50%  '
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('D41').Value = "'2.40"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = "'173.47"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.34%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = "'0.000309"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('D48').Value = "'49.53"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('D49').Value = "'424.44"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('D51').Value = "'8.54"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.01%  '

Write-Host "Applied cryptos update"
